$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Labels")

# Insert a new row for the "DISCHARGE" label substitution (reuses the \f
# escape code, right below the existing FILL row) and push the rows that
# used to follow (RELEASE, HEATING, COOLING, FLAP, CONTROL, ...) down by one.
$ws.Rows.Item(28).Insert()
$ws.Range("A28").Value = "\\f"
$ws.Range("B28").Value = "DISCHARGE"

# Match the author's new cursor/selection position on the sheet.
$ws.Activate()
$ws.Range("B27").Select() | Out-Null
